# edit.ps1
# Updates the cryptocurrency price/volume table on the active sheet.
#
# The source data refreshed: most rows keep the same coin but get new
# Price (D) / Volume(1h) (E) figures; additionally a new row for
# "Frax" was inserted at row 34, which pushes every following coin
# down by one row and drops the former last row ("Aave", row 51) off
# the bottom of the (unchanged) A1:E51 sheet range. Column A (the
# original index numbers) is left untouched.
#
# Many of the Price values look like plain numbers (e.g. "1.001",
# "244.19") even though the source file stores them as literal text
# (some, like "26.491.44", have two decimal points and are not valid
# numbers at all). Writing such a string straight into a cell's
# .Value causes Excel to auto-convert it into a numeric value, which
# would corrupt values such as "26.491.44" and strip significant
# trailing zeros (e.g. "30.30" -> 30.3). To avoid this we temporarily
# mark the cell as Text ("@") before assigning the value, then revert
# the cell's format back to General/Normal so no stray formatting is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Each entry: row, Coin (B), Link (C), Price (D), Volume(1h) (E)
$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "26.491.44", "  -0.11%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.731.88", "  +0.23%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.001", "  +0.21%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "244.19", "  -0.37%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  +0.15%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4907", "  +2.10%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2627", "  -1.59%  "),
    @(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06187", "  -0.55%  "),
    @(10, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.733.62", "  +0.40%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07014", "  -1.89%  "),
    @(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "15.49", "  -1.04%  "),
    @(13, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.555", "  +0.65%  "),
    @(14, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6011", "  -2.53%  "),
    @(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "77.49", "  +0.42%  "),
    @(16, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  +0.14%  "),
    @(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "26.511.47", "  -0.06%  "),
    @(18, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.002", "  +0.18%  "),
    @(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007089", "  +2.28%  "),
    @(20, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "11.42", "  -2.10%  "),
    @(21, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.958.85", "  +0.60%  "),
    @(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.488", "  -0.90%  "),
    @(23, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.610", "  -3.81%  "),
    @(24, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.186", "  -1.86%  "),
    @(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "138.83", "  +1.65%  "),
    @(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "15.23", "  -0.72%  "),
    @(27, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.415", "  +0.76%  "),
    @(28, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "106.62", "  -0.12%  "),
    @(29, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.716", "  -4.38%  "),
    @(30, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.986", "  +0.24%  "),
    @(31, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.07962", "  -0.81%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.701", "  -0.19%  "),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04529", "  -0.79%  "),
    @(34, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "1.001", "  +0.14%  "),
    @(35, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.615", "  +0.03%  "),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.001", "  +0.91%  "),
    @(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.6256", "  -1.84%  "),
    @(38, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9048", "  -2.82%  "),
    @(39, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.999", "  -4.67%  "),
    @(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.408", "  -0.58%  "),
    @(41, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.003", "  -0.32%  "),
    @(42, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01488", "  -1.00%  "),
    @(43, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "100.52", "  -4.12%  "),
    @(44, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.469", "  -2.19%  "),
    @(45, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3869", "  -0.95%  "),
    @(46, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.653", "  -3.68%  "),
    @(47, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1157", "  -2.32%  "),
    @(48, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05364", "  +0.63%  "),
    @(49, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "30.30", "  -1.87%  "),
    @(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "7.665", "  -2.20%  "),
    @(51, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.249", "  -1.63%  ")
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    Set-TextValue $ws.Cells.Item($r, 4) $entry[3]
    Set-TextValue $ws.Cells.Item($r, 5) $entry[4]
}
